$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "29.210.98"
Set-TextCell "E2" "  +0.23%  "
Set-TextCell "D3" "1.832.84"
Set-TextCell "E3" "  -0.41%  "
Set-TextCell "D4" "0.9986"
Set-TextCell "E4" "  -0.25%  "
Set-TextCell "D5" "242.42"
Set-TextCell "E5" "  -0.66%  "
Set-TextCell "D6" "0.6234"
Set-TextCell "E6" "  -0.48%  "
Set-TextCell "D7" "0.9996"
Set-TextCell "E7" "  -0.20%  "
Set-TextCell "D8" "0.07378"
Set-TextCell "E8" "  -1.79%  "
Set-TextCell "D9" "0.2908"
Set-TextCell "E9" "  -1.08%  "
Set-TextCell "D10" "23.14"
Set-TextCell "E10" "  -0.73%  "
Set-TextCell "D11" "0.07673"
Set-TextCell "E11" "  -0.51%  "
Set-TextCell "D12" "1.829.17"
Set-TextCell "E12" "  -1.51%  "
Set-TextCell "D13" "4.965"
Set-TextCell "E13" "  -1.16%  "
Set-TextCell "D14" "0.6670"
Set-TextCell "D15" "82.56"
Set-TextCell "E15" "  -0.61%  "
Set-TextCell "D16" "0.000008966"
Set-TextCell "E16" "  -3.28%  "
Set-TextCell "D17" "5.871"
Set-TextCell "E17" "  -1.85%  "
Set-TextCell "D18" "29.165.93"
Set-TextCell "E18" "  +0.04%  "
Set-TextCell "D19" "2.079.56"
Set-TextCell "E19" "  -2.12%  "
Set-TextCell "D20" "235.97"
Set-TextCell "E20" "  +2.31%  "
Set-TextCell "D22" "0.9995"
Set-TextCell "E22" "  -0.24%  "
Set-TextCell "D23" "7.398"
Set-TextCell "E23" "  +2.89%  "
Set-TextCell "D24" "0.9999"
Set-TextCell "E24" "  -0.22%  "
Set-TextCell "D25" "158.19"
Set-TextCell "E25" "  -1.42%  "
Set-TextCell "D26" "0.1411"
Set-TextCell "E26" "  +1.44%  "
Set-TextCell "E27" "  -0.43%  "
Set-TextCell "D28" "17.65"
Set-TextCell "E28" "  -1.41%  "
Set-TextCell "D29" "1.484"
Set-TextCell "E29" "  -1.07%  "
Set-TextCell "D30" "0.05798"
Set-TextCell "E30" "  +4.32%  "
Set-TextCell "D31" "4.096"
Set-TextCell "D32" "4.089"
Set-TextCell "E32" "  -2.54%  "
Set-TextCell "D33" "1.205"
Set-TextCell "E33" "  -0.41%  "
Set-TextCell "E34" "  +0.58%  "
Set-TextCell "D35" "0.7320"
Set-TextCell "E35" "  -2.39%  "
Set-TextCell "E36" "  -0.62%  "
Set-TextCell "D37" "2.607"
Set-TextCell "E37" "  -1.99%  "
Set-TextCell "D38" "2.841"
Set-TextCell "E38" "  +2.42%  "
Set-TextCell "D39" "1.226.78"
Set-TextCell "E39" "  -0.36%  "
Set-TextCell "E40" "  -1.99%  "
Set-TextCell "B41" "TrustWalletToken"
Set-TextCell "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D41" "0.9168"
Set-TextCell "E41" "  +1.69%  "
Set-TextCell "B42" "FraxShare"
Set-TextCell "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D42" "6.262"
Set-TextCell "E42" "  -4.62%  "
Set-TextCell "D43" "1.000"
Set-TextCell "E43" "  -0.05%  "
Set-TextCell "D44" "101.90"
Set-TextCell "E44" "  -0.32%  "
Set-TextCell "D45" "1.981.15"
Set-TextCell "E45" "  -1.76%  "
Set-TextCell "D46" "65.08"
Set-TextCell "E46" "  -1.99%  "
Set-TextCell "D47" "0.5041"
Set-TextCell "E47" "  -1.21%  "
Set-TextCell "E48" "  -4.62%  "
Set-TextCell "D49" "0.4026"
Set-TextCell "D50" "9.113"
Set-TextCell "E51" "  +2.57%  "
